$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, pushing existing rows 193:226 down to 194:227
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new weekly price record
$ws.Range("A193").Value = 3
$ws.Range("B193").Value = 'Femacal de La Calera'
$ws.Range("C193").Value = 'Coquimbo'
$ws.Range("D193").Value = 44476
$ws.Range("E193").Value = 5
$ws.Range("F193").Value = 100112003
$ws.Range("G193").Value = 'Ajo'
$ws.Range("H193").Value = 'Chino'
$ws.Range("I193").Value = 'Primera'
$ws.Range("J193").Value = 110
$ws.Range("K193").Value = 16500
$ws.Range("L193").Value = 17000
$ws.Range("M193").Value = 16727
$ws.Range("N193").Value = '$/caja 10 kilos'
$ws.Range("O193").Value = 'China'
$ws.Range("P193").Value = 1673
$ws.Range("Q193").Value = 10
$ws.Range("R193").Value = 'Hortaliza'
